$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- 1) Add the two new trailing columns (I, J) with headers first so the ---
# --- new shared strings land in the same order as the source edit.      ---
$ws.Range("I1").Value = "LontitudeMain"
$ws.Range("J1").Value = "LatitudeMain"
# "Water Education Colorado" is row 7 before the insert below shifts it down
# to row 8, so the new lon/lat-main pair must be written at I7/J7 now.
$ws.Range("I7").Value = -104.972905
$ws.Range("J7").Value = 39.742072999999998

$ws.Columns("I").ColumnWidth = 14.6666666666667
$ws.Columns("J").ColumnWidth = 19

# --- 2) Insert the new "Environmental Learning Center" row at row 3,    ---
# --- pushing the existing rows 3-8 down to 4-9.                        ---
$ws.Rows("3:3").Insert()

$ws.Range("A3").Value = "Environmental Learning Center"
$ws.Range("B3").Value = "University"
$ws.Range("C3").Value = "University research with community access."
$ws.Range("D3").Value = "River habitat education, wildlife."
$ws.Range("E3").Value = "https://warnercnr.colostate.edu/elc/"
$ws.Range("F3").Value = "Yes"
$ws.Range("G3").Value = -105.019469
$ws.Range("H3").Value = 40.555596999999999

# --- 3) The "InBasin" column (F) loses its Hyperlink styling. ---
$ws.Range("F2:F9").Style = "Normal"

# --- 4) Rebuild the hyperlinks collection: ranges shifted down one row for ---
# --- the original 7 links, plus a brand-new one for the inserted row.     ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E4"), "https://www.northernwater.org/AboutUs/WaterEducation.aspx")
$ws.Hyperlinks.Add($ws.Range("E5"), "https://www.poudreheritage.org/")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://watercenter.colostate.edu/")
$ws.Hyperlinks.Add($ws.Range("E6"), "https://poudrelearningcenter.org/")
$ws.Hyperlinks.Add($ws.Range("E7"), "https://southplattebasin.com/")
$ws.Hyperlinks.Add($ws.Range("E8"), "https://www.watereducationcolorado.org/")
$ws.Hyperlinks.Add($ws.Range("E9"), "http://openwaterfoundation.org/")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://warnercnr.colostate.edu/elc/")

# Hyperlinks.Add() re-stamps a fresh "hyperlink-ish" style on the target
# cell; put the standard Hyperlink style back so the column keeps its
# original formatting.
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("E3").Style = "Hyperlink"
$ws.Range("E4").Style = "Hyperlink"
$ws.Range("E5").Style = "Hyperlink"
$ws.Range("E6").Style = "Hyperlink"
$ws.Range("E7").Style = "Hyperlink"
$ws.Range("E8").Style = "Hyperlink"
$ws.Range("E9").Style = "Hyperlink"
